$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ventas")

# --- Insert the first new data row (42400 / 9) above the current row 2 ---
$ws.Rows("2:2").Insert()
$ws.Range("A2:B2").ClearFormats()
$ws.Range("A2").Value = 42400
$ws.Range("B2").Value = 9
$ws.Range("A2").NumberFormat = "yyyy/mm/dd"

# --- Insert the second new data row (42401 / 2) above what is now row 8 ---
# (the row that still holds the old last record, 42388 / 3)
$ws.Rows("8:8").Insert()
$ws.Range("A8:B8").ClearFormats()
$ws.Range("A8").Value = 42401
$ws.Range("B8").Value = 2
$ws.Range("A8").NumberFormat = "yyyy/mm/dd"

# --- Update the line chart's source ranges so they cover the new rows ---
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(,Ventas!`$A`$2:`$A`$9,Ventas!`$B`$2:`$B`$9,1)"
